# Scheduled runner: refresh cached Universalis market-price figures
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns)
# for the specific leve rows whose market data changed since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 5009
$ws.Range("I38").Value = 5009
$ws.Range("K38").Value = 15027
$ws.Range("M38").Value = -14655

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 356.66666
$ws.Range("I58").Value = 356.66666
$ws.Range("K58").Value = 1069.99998
$ws.Range("M58").Value = -919.9999800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 13994.75
$ws.Range("I69").Value = 7999
$ws.Range("J69").Value = 15993.333
$ws.Range("K69").Value = 23997
$ws.Range("L69").Value = 47979.999
$ws.Range("M69").Value = -23123
$ws.Range("N69").Value = -49727.999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 13994.75
$ws.Range("I72").Value = 7999
$ws.Range("J72").Value = 15993.333
$ws.Range("K72").Value = 71991
$ws.Range("L72").Value = 143939.997
$ws.Range("M72").Value = -67623
$ws.Range("N72").Value = -152675.997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1977.7778
$ws.Range("J112").Value = 1977.7778
$ws.Range("L112").Value = 5933.3334
$ws.Range("N112").Value = -8149.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 997.6
$ws.Range("J121").Value = 997.6
$ws.Range("L121").Value = 2992.8
$ws.Range("N121").Value = -6486.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2416.9412
$ws.Range("I135").Value = 2149.5715
$ws.Range("K135").Value = 19346.1435
$ws.Range("M135").Value = -16811.1435

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2692.3845
$ws.Range("I137").Value = 2077.5625
$ws.Range("J137").Value = 3676.1
$ws.Range("K137").Value = 6232.6875
$ws.Range("L137").Value = 11028.3
$ws.Range("M137").Value = -3682.6875
$ws.Range("N137").Value = -16128.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 10668.138
$ws.Range("J138").Value = 10693.306
$ws.Range("L138").Value = 32079.918
$ws.Range("N138").Value = -42359.91800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 4500
$ws.Range("I26").Value = 4500
$ws.Range("K26").Value = 4500
$ws.Range("M26").Value = -4170

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17234.393
$ws.Range("I32").Value = 13576.5
$ws.Range("J32").Value = 26013.334
$ws.Range("K32").Value = 13576.5
$ws.Range("L32").Value = 26013.334
$ws.Range("M32").Value = -13289.5
$ws.Range("N32").Value = -26587.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8551.041999999999
$ws.Range("I61").Value = 6366.2354
$ws.Range("K61").Value = 6366.2354
$ws.Range("M61").Value = -6154.2354

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3258.8125
$ws.Range("I74").Value = 2374.3572
$ws.Range("K74").Value = 2374.3572
$ws.Range("M74").Value = -1500.3572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3258.8125
$ws.Range("I77").Value = 2374.3572
$ws.Range("K77").Value = 11871.786
$ws.Range("M77").Value = -7503.786

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 46750
$ws.Range("J113").Value = 46750
$ws.Range("L113").Value = 46750
$ws.Range("N113").Value = -55428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4380
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4014.2766
$ws.Range("I132").Value = 3259.6711
$ws.Range("J132").Value = 7200.3887
$ws.Range("K132").Value = 9779.013300000001
$ws.Range("L132").Value = 21601.1661
$ws.Range("M132").Value = -7249.013300000001
$ws.Range("N132").Value = -26661.1661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8551.041999999999
$ws.Range("I136").Value = 6366.2354
$ws.Range("K136").Value = 19098.7062
$ws.Range("M136").Value = -16548.7062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1249.3334
$ws.Range("I99").Value = 1210.8462
$ws.Range("K99").Value = 1210.8462
$ws.Range("M99").Value = 287.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 13185.143
$ws.Range("I107").Value = 13028.765
$ws.Range("K107").Value = 13028.765
$ws.Range("M107").Value = -11108.765

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 18516.096
$ws.Range("I134").Value = 4359.9683
$ws.Range("K134").Value = 13079.9049
$ws.Range("M134").Value = -10544.9049

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4800.8
$ws.Range("I99").Value = 4051.0833
$ws.Range("K99").Value = 4051.0833
$ws.Range("M99").Value = -2553.0833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4800.8
$ws.Range("I126").Value = 4051.0833
$ws.Range("K126").Value = 12153.2499
$ws.Range("M126").Value = -9683.249899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3598.9375
$ws.Range("I132").Value = 3121.7693
$ws.Range("J132").Value = 5666.6665
$ws.Range("K132").Value = 9365.3079
$ws.Range("L132").Value = 16999.9995
$ws.Range("M132").Value = -6835.3079
$ws.Range("N132").Value = -22059.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 438098.2
$ws.Range("I134").Value = 3517.0476
$ws.Range("K134").Value = 10551.1428
$ws.Range("M134").Value = -8016.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 565714
$ws.Range("I141").Value = 60000
$ws.Range("J141").Value = 767999.6
$ws.Range("K141").Value = 60000
$ws.Range("L141").Value = 767999.6
$ws.Range("M141").Value = -54820
$ws.Range("N141").Value = -778359.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 9002601
$ws.Range("I32").Value = 10002900
$ws.Range("K32").Value = 30008700
$ws.Range("M32").Value = -30008417

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 11547.72
$ws.Range("I39").Value = 4814
$ws.Range("K39").Value = 14442
$ws.Range("M39").Value = -14148

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 6912.229
$ws.Range("I139").Value = 3237.625
$ws.Range("K139").Value = 9712.875
$ws.Range("M139").Value = -4572.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 24094.666
$ws.Range("J26").Value = 25493.846
$ws.Range("L26").Value = 25493.846
$ws.Range("N26").Value = -26053.846

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 23081.766
$ws.Range("I43").Value = 10090
$ws.Range("K43").Value = 10090
$ws.Range("M43").Value = -9939

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 24094.666
$ws.Range("J50").Value = 25493.846
$ws.Range("L50").Value = 25493.846
$ws.Range("N50").Value = -26489.846

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 45461220
$ws.Range("I70").Value = 6313.615
$ws.Range("K70").Value = 6313.615
$ws.Range("M70").Value = -6043.615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 45461220
$ws.Range("I73").Value = 6313.615
$ws.Range("K73").Value = 6313.615
$ws.Range("M73").Value = -5377.615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5001250
$ws.Range("I80").Value = 3334998.2
$ws.Range("K80").Value = 3334998.2
$ws.Range("M80").Value = -3334000.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5001250
$ws.Range("I83").Value = 3334998.2
$ws.Range("K83").Value = 16674991
$ws.Range("M83").Value = -16669999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2887.5
$ws.Range("I46").Value = 2516.6667
$ws.Range("K46").Value = 2516.6667
$ws.Range("M46").Value = -2328.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 49799.8
$ws.Range("J133").Value = 49799.8
$ws.Range("L133").Value = 49799.8
$ws.Range("N133").Value = -54859.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9290.177
$ws.Range("I136").Value = 7350.905
$ws.Range("J136").Value = 12422.846
$ws.Range("K136").Value = 22052.715
$ws.Range("L136").Value = 37268.538
$ws.Range("M136").Value = -19502.715
$ws.Range("N136").Value = -42368.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 50409
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 50409
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 50409
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -59587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 27961.07
$ws.Range("I132").Value = 4391.5
$ws.Range("K132").Value = 13174.5
$ws.Range("M132").Value = -10644.5
